$wb = $excel.ActiveWorkbook

# ---- Worksheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 76
$ws.Cells.Item(76, 8).Value = 3370
$ws.Cells.Item(76, 9).Value = 3370
$ws.Cells.Item(76, 11).Value = 3370
$ws.Cells.Item(76, 13).Value = -3055
# Row 79
$ws.Cells.Item(79, 8).Value = 3370
$ws.Cells.Item(79, 9).Value = 3370
$ws.Cells.Item(79, 11).Value = 3370
$ws.Cells.Item(79, 13).Value = -2278
# Row 137
$ws.Cells.Item(137, 8).Value = 1935.2766
$ws.Cells.Item(137, 9).Value = 1461.6875
$ws.Cells.Item(137, 10).Value = 2945.6
$ws.Cells.Item(137, 11).Value = 4385.0625
$ws.Cells.Item(137, 12).Value = 8836.799999999999
$ws.Cells.Item(137, 13).Value = -1835.0625
$ws.Cells.Item(137, 14).Value = -13936.8
# Row 138
$ws.Cells.Item(138, 8).Value = 3638.4614
$ws.Cells.Item(138, 9).Value = 2833.9285
$ws.Cells.Item(138, 10).Value = 4577.0835
$ws.Cells.Item(138, 11).Value = 8501.7855
$ws.Cells.Item(138, 12).Value = 13731.2505
$ws.Cells.Item(138, 13).Value = -3361.7855
$ws.Cells.Item(138, 14).Value = -24011.2505
# Row 141
$ws.Cells.Item(141, 8).Value = 1112872.1
$ws.Cells.Item(141, 9).Value = 1155290.2
$ws.Cells.Item(141, 11).Value = 3465870.6
$ws.Cells.Item(141, 13).Value = -3460690.6

# ---- Worksheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 74
$ws.Cells.Item(74, 8).Value = 86534.38
$ws.Cells.Item(74, 9).Value = 143758.28
$ws.Cells.Item(74, 10).Value = 19773.166
$ws.Cells.Item(74, 11).Value = 143758.28
$ws.Cells.Item(74, 12).Value = 19773.166
$ws.Cells.Item(74, 13).Value = -142884.28
$ws.Cells.Item(74, 14).Value = -21521.166
# Row 77
$ws.Cells.Item(77, 8).Value = 86534.38
$ws.Cells.Item(77, 9).Value = 143758.28
$ws.Cells.Item(77, 10).Value = 19773.166
$ws.Cells.Item(77, 11).Value = 718791.4
$ws.Cells.Item(77, 12).Value = 98865.83
$ws.Cells.Item(77, 13).Value = -714423.4
$ws.Cells.Item(77, 14).Value = -107601.83

# ---- Worksheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Cells.Item(134, 8).Value = 550542.6
$ws.Cells.Item(134, 9).Value = 664402.6
$ws.Cells.Item(134, 10).Value = 7518
$ws.Cells.Item(134, 11).Value = 1993207.8
$ws.Cells.Item(134, 12).Value = 22554
$ws.Cells.Item(134, 13).Value = -1990672.8
$ws.Cells.Item(134, 14).Value = -27624

# ---- Worksheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 744.3542
$ws.Cells.Item(31, 9).Value = 702.5054
$ws.Cells.Item(31, 10).Value = 2041.6666
$ws.Cells.Item(31, 11).Value = 702.5054
$ws.Cells.Item(31, 12).Value = 2041.6666
$ws.Cells.Item(31, 13).Value = -407.5054
$ws.Cells.Item(31, 14).Value = -2631.6666
# Row 34
$ws.Cells.Item(34, 8).Value = 744.3542
$ws.Cells.Item(34, 9).Value = 702.5054
$ws.Cells.Item(34, 10).Value = 2041.6666
$ws.Cells.Item(34, 11).Value = 702.5054
$ws.Cells.Item(34, 12).Value = 2041.6666
$ws.Cells.Item(34, 13).Value = -500.5054
$ws.Cells.Item(34, 14).Value = -2445.6666
# Row 132
$ws.Cells.Item(132, 8).Value = 2647857.5
$ws.Cells.Item(132, 9).Value = 1839.2
$ws.Cells.Item(132, 11).Value = 5517.6
$ws.Cells.Item(132, 13).Value = -2987.6

# ---- Worksheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 122
$ws.Cells.Item(122, 8).Value = 28876.871
$ws.Cells.Item(122, 9).Value = 34834.375
$ws.Cells.Item(122, 10).Value = 1642.5714
$ws.Cells.Item(122, 11).Value = 313509.375
$ws.Cells.Item(122, 12).Value = 14783.1426
$ws.Cells.Item(122, 13).Value = -311059.375
$ws.Cells.Item(122, 14).Value = -19683.1426
# Row 128
$ws.Cells.Item(128, 8).Value = 55666.668
$ws.Cells.Item(128, 9).Value = 55666.668
$ws.Cells.Item(128, 11).Value = 167000.004
$ws.Cells.Item(128, 13).Value = -162020.004
# Row 134
$ws.Cells.Item(134, 8).Value = 2676.818
$ws.Cells.Item(134, 9).Value = 2566.1904
$ws.Cells.Item(134, 10).Value = 5000
$ws.Cells.Item(134, 11).Value = 7698.5712
$ws.Cells.Item(134, 12).Value = 15000
$ws.Cells.Item(134, 13).Value = -2628.5712
$ws.Cells.Item(134, 14).Value = -25140
# Row 139
$ws.Cells.Item(139, 8).Value = 2277.4119
$ws.Cells.Item(139, 9).Value = 1647.7333
$ws.Cells.Item(139, 10).Value = 7000
$ws.Cells.Item(139, 11).Value = 4943.199900000001
$ws.Cells.Item(139, 12).Value = 21000
$ws.Cells.Item(139, 13).Value = 196.8000999999995
$ws.Cells.Item(139, 14).Value = -31280

# ---- Worksheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 80
$ws.Cells.Item(80, 8).Value = 20063
$ws.Cells.Item(80, 9).Value = 19998
$ws.Cells.Item(80, 10).Value = 20128
$ws.Cells.Item(80, 11).Value = 19998
$ws.Cells.Item(80, 12).Value = 20128
$ws.Cells.Item(80, 13).Value = -18875
$ws.Cells.Item(80, 14).Value = -22374
# Row 81
$ws.Cells.Item(81, 8).Value = 42000
$ws.Cells.Item(81, 9).Value = 20000
$ws.Cells.Item(81, 10).Value = 47500
$ws.Cells.Item(81, 11).Value = 20000
$ws.Cells.Item(81, 12).Value = 47500
$ws.Cells.Item(81, 13).Value = -19002
$ws.Cells.Item(81, 14).Value = -49496
# Row 83
$ws.Cells.Item(83, 8).Value = 20063
$ws.Cells.Item(83, 9).Value = 19998
$ws.Cells.Item(83, 10).Value = 20128
$ws.Cells.Item(83, 11).Value = 59994
$ws.Cells.Item(83, 12).Value = 60384
$ws.Cells.Item(83, 13).Value = -54378
$ws.Cells.Item(83, 14).Value = -71616
# Row 84
$ws.Cells.Item(84, 8).Value = 42000
$ws.Cells.Item(84, 9).Value = 20000
$ws.Cells.Item(84, 10).Value = 47500
$ws.Cells.Item(84, 11).Value = 60000
$ws.Cells.Item(84, 12).Value = 142500
$ws.Cells.Item(84, 13).Value = -55008
$ws.Cells.Item(84, 14).Value = -152484
# Row 92
$ws.Cells.Item(92, 8).Value = 29798.428
$ws.Cells.Item(92, 10).Value = 29798.428
$ws.Cells.Item(92, 12).Value = 29798.428
$ws.Cells.Item(92, 14).Value = -34790.428
# Row 98
$ws.Cells.Item(98, 8).Value = 29118.334
$ws.Cells.Item(98, 10).Value = 29118.334
$ws.Cells.Item(98, 12).Value = 29118.334
$ws.Cells.Item(98, 14).Value = -35108.334
# Row 104
$ws.Cells.Item(104, 8).Value = 27500
$ws.Cells.Item(104, 10).Value = 27500
$ws.Cells.Item(104, 12).Value = 27500
$ws.Cells.Item(104, 14).Value = -34488
# Row 108
$ws.Cells.Item(108, 8).Value = 0
$ws.Cells.Item(108, 10).Value = 0
$ws.Cells.Item(108, 12).Value = 0
$ws.Cells.Item(108, 14).ClearContents()
# Row 110
$ws.Cells.Item(110, 8).Value = 0
$ws.Cells.Item(110, 10).Value = 0
$ws.Cells.Item(110, 12).Value = 0
$ws.Cells.Item(110, 14).ClearContents()
# Row 114
$ws.Cells.Item(114, 8).Value = 0
$ws.Cells.Item(114, 10).Value = 0
$ws.Cells.Item(114, 12).Value = 0
$ws.Cells.Item(114, 14).ClearContents()
# Row 116
$ws.Cells.Item(116, 8).Value = 0
$ws.Cells.Item(116, 10).Value = 0
$ws.Cells.Item(116, 12).Value = 0
$ws.Cells.Item(116, 14).ClearContents()
# Row 120
$ws.Cells.Item(120, 8).Value = 31142.857
$ws.Cells.Item(120, 10).Value = 31142.857
$ws.Cells.Item(120, 12).Value = 31142.857
$ws.Cells.Item(120, 14).Value = -40818.857
# Row 122
$ws.Cells.Item(122, 8).Value = 2637.5
$ws.Cells.Item(122, 9).Value = 2366.6667
$ws.Cells.Item(122, 10).Value = 3450
$ws.Cells.Item(122, 11).Value = 7100.000100000001
$ws.Cells.Item(122, 12).Value = 10350
$ws.Cells.Item(122, 13).Value = -4650.000100000001
$ws.Cells.Item(122, 14).Value = -15250
# Row 123
$ws.Cells.Item(123, 8).Value = 24285.715
$ws.Cells.Item(123, 10).Value = 24285.715
$ws.Cells.Item(123, 12).Value = 24285.715
$ws.Cells.Item(123, 14).Value = -34085.715
# Row 124
$ws.Cells.Item(124, 8).Value = 30000
$ws.Cells.Item(124, 10).Value = 30000
$ws.Cells.Item(124, 12).Value = 30000
$ws.Cells.Item(124, 14).Value = -39820
# Row 125
$ws.Cells.Item(125, 8).Value = 20000
$ws.Cells.Item(125, 10).Value = 20000
$ws.Cells.Item(125, 12).Value = 20000
$ws.Cells.Item(125, 14).Value = -29840
# Row 127
$ws.Cells.Item(127, 8).Value = 35101.668
$ws.Cells.Item(127, 10).Value = 35101.668
$ws.Cells.Item(127, 12).Value = 35101.668
$ws.Cells.Item(127, 14).Value = -45021.668
# Row 128
$ws.Cells.Item(128, 8).Value = 37142.855
$ws.Cells.Item(128, 10).Value = 37142.855
$ws.Cells.Item(128, 12).Value = 37142.855
$ws.Cells.Item(128, 14).Value = -47102.855
# Row 129
$ws.Cells.Item(129, 8).Value = 0
$ws.Cells.Item(129, 10).Value = 0
$ws.Cells.Item(129, 12).Value = 0
$ws.Cells.Item(129, 14).ClearContents()

# ---- Worksheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Cells.Item(81, 8).Value = 1777.8
$ws.Cells.Item(81, 9).Value = 2267
$ws.Cells.Item(81, 10).Value = 1044
$ws.Cells.Item(81, 11).Value = 4534
$ws.Cells.Item(81, 12).Value = 2088
$ws.Cells.Item(81, 13).Value = -3473
$ws.Cells.Item(81, 14).Value = -4210
# Row 84
$ws.Cells.Item(84, 8).Value = 1777.8
$ws.Cells.Item(84, 9).Value = 2267
$ws.Cells.Item(84, 10).Value = 1044
$ws.Cells.Item(84, 11).Value = 22670
$ws.Cells.Item(84, 12).Value = 10440
$ws.Cells.Item(84, 13).Value = -17366
$ws.Cells.Item(84, 14).Value = -21048
# Row 113
$ws.Cells.Item(113, 8).Value = 476.23077
$ws.Cells.Item(113, 9).Value = 528
$ws.Cells.Item(113, 10).Value = 443.875
$ws.Cells.Item(113, 11).Value = 1584
$ws.Cells.Item(113, 12).Value = 1331.625
$ws.Cells.Item(113, 13).Value = 586
$ws.Cells.Item(113, 14).Value = -5671.625
# Row 136
$ws.Cells.Item(136, 8).Value = 5436.12
$ws.Cells.Item(136, 9).Value = 942
$ws.Cells.Item(136, 10).Value = 19667.5
$ws.Cells.Item(136, 11).Value = 2826
$ws.Cells.Item(136, 12).Value = 59002.5
$ws.Cells.Item(136, 13).Value = -276
$ws.Cells.Item(136, 14).Value = -64102.5
